# Update countries & provincias Spain
#
# Refreshes the "Pais" sheet with a newer COVID-19 data pull:
#   - bumps the "last updated" timestamp
#   - updates the running totals (Casos totales/Nuevos casos/Casos activos/
#     Recuperados/Casos criticos/Muertes hoy/Muertes) for the countries whose
#     counts moved
#   - a few countries overtook their neighbour in the ranking, so the two
#     rows in that pair swap which country name they show (their stats are
#     rewritten to match)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 02:31"

# --- Straight numeric refreshes (country/ranking unchanged) ----------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2233957
$ws.Range("C4").Value = 25557
$ws.Range("D4").Value = 912741
$ws.Range("E4").Value = 1201275
$ws.Range("G4").Value = 809
$ws.Range("H4").Value = 119941

# Row 5 - Brasil
$ws.Range("D5").Value = 503507
$ws.Range("E5").Value = 410137

# Row 37 - Argentina
$ws.Range("B37").Value = 35552
$ws.Range("C37").Value = 1393
$ws.Range("E37").Value = 24127
$ws.Range("G37").Value = 35
$ws.Range("H37").Value = 913

# Row 39 - Suiza
$ws.Range("B39").Value = 31187
$ws.Range("C39").Value = 33
$ws.Range("E39").Value = 331

# Row 48 - Panama
$ws.Range("B48").Value = 22597
$ws.Range("C48").Value = 635
$ws.Range("E48").Value = 8353
$ws.Range("G48").Value = 13
$ws.Range("H48").Value = 470

# Row 84 - Gabon
$ws.Range("B84").Value = 4229
$ws.Range("C84").Value = 115
$ws.Range("D84").Value = 1505
$ws.Range("E84").Value = 2694
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 30

# Row 87 - El Salvador
$ws.Range("E87").Value = 1850
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 79

# Row 126 - Niger
$ws.Range("B126").Value = 1020
$ws.Range("C126").Value = 4
$ws.Range("D126").Value = 893
$ws.Range("E126").Value = 60
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 67

# Row 143 - Ruanda
$ws.Range("B143").Value = 639
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 347
$ws.Range("E143").Value = 290

# Row 149 - Togo
$ws.Range("B149").Value = 544
$ws.Range("C149").Value = 7
$ws.Range("D149").Value = 353
$ws.Range("E149").Value = 178

# Row 161 - Surinam
$ws.Range("B161").Value = 261
$ws.Range("C161").Value = 25
$ws.Range("E161").Value = 207

# --- Ranking swaps: country name + stats trade places between the two
#     rows of each pair ------------------------------------------------

# Rows 109/110: Sudan del Sur overtakes Lituania
$ws.Range("A109").Value = "Sudan del Sur"
$ws.Range("B109").Value = 1813
$ws.Range("C109").Value = 37
$ws.Range("D109").Value = 89
$ws.Range("E109").Value = 1693
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 31
$ws.Range("A110").Value = "Lituania"
$ws.Range("B110").Value = 1778
$ws.Range("C110").Value = 2
$ws.Range("D110").Value = 1447
$ws.Range("E110").Value = 255
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 76

# Rows 145/146: Malaui overtakes Benin
$ws.Range("A145").Value = "Malaui"
$ws.Range("B145").Value = 572
$ws.Range("C145").Value = 8
$ws.Range("D145").Value = 73
$ws.Range("E145").Value = 493
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 6
$ws.Range("A146").Value = "Benin"
$ws.Range("B146").Value = 572
$ws.Range("C146").Value = 40
$ws.Range("D146").Value = 237
$ws.Range("E146").Value = 326
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 9

# Rows 152/153: Libia overtakes Reunion
$ws.Range("A152").Value = "Libia"
$ws.Range("B152").Value = 500
$ws.Range("C152").Value = 16
$ws.Range("D152").Value = 78
$ws.Range("E152").Value = 412
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 10
$ws.Range("A153").Value = "Reunion"
$ws.Range("B153").Value = 497
$ws.Range("C153").Value = 2
$ws.Range("D153").Value = 460
$ws.Range("E153").Value = 36
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 1

# Rows 206/207: Islas Malvinas overtakes Groenlandia
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("B206").Value = 13
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 13
$ws.Range("E206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("B207").Value = 13
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 13
$ws.Range("E207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

# Rows 208/209: Santa Sede overtakes Islas Turcas y Caicos
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 12
$ws.Range("E208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("B209").Value = 12
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("E209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1
